# Rotate observation-record rows 45-48 and 57-59 (content shifts down by one row,
# with the last row of each block wrapping around to become the first).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 45 ---
$ws.Range("A45").Value = 131031771
$ws.Range("B45").Value = 79243
$ws.Range("E45").Value = 6425
$ws.Range("F45").Value = 'Garnlav'
$ws.Range("G45").Value = 'Alectoria sarmentosa'
$ws.Range("H45").Value = '(Ach.) Ach.'
$ws.Range("M45").ClearContents()
$ws.Range("P45").Value = 'Stacktjärnen, Stacktjärnen, Jmt'
$ws.Range("Q45").Value = 468053
$ws.Range("R45").Value = 7017383
$ws.Range("Z45").Value = '14:50'
$ws.Range("AB45").Value = '14:50'
$ws.Range("AC45").ClearContents()

# --- Row 46 ---
$ws.Range("A46").Value = 131030403
$ws.Range("Q46").Value = 468347
$ws.Range("R46").Value = 7017234
$ws.Range("S46").Value = 10
$ws.Range("Z46").Value = '13:49'
$ws.Range("AB46").Value = '13:49'
$ws.Range("AC46").Value = 'Barkfläkta granar i området'

# --- Row 47 ---
$ws.Range("A47").Value = 131030207
$ws.Range("M47").Value = 'äldre spår'
$ws.Range("Q47").Value = 468282
$ws.Range("R47").Value = 7017338
$ws.Range("S47").Value = 5
$ws.Range("Z47").Value = '13:33'
$ws.Range("AB47").Value = '13:33'
$ws.Range("AC47").Value = 'Ringhack'

# --- Row 48 ---
$ws.Range("A48").Value = 131030656
$ws.Range("B48").Value = 57884
$ws.Range("E48").Value = 100109
$ws.Range("F48").Value = 'Tretåig hackspett'
$ws.Range("G48").Value = 'Picoides tridactylus'
$ws.Range("H48").Value = '(Linnaeus, 1758)'
$ws.Range("M48").Value = 'färska spår'
$ws.Range("P48").Value = 'Stacktjärnen, Rödön, Stacktjärnen, Rödön, Jmt'
$ws.Range("Q48").Value = 468332
$ws.Range("R48").Value = 7017164
$ws.Range("Z48").Value = '13:59'
$ws.Range("AB48").Value = '13:59'
$ws.Range("AC48").Value = 'Gott om barkfläkta granar i området'

# --- Row 57 ---
$ws.Range("A57").Value = 131031498
$ws.Range("B57").Value = 57884
$ws.Range("E57").Value = 100109
$ws.Range("F57").Value = 'Tretåig hackspett'
$ws.Range("G57").Value = 'Picoides tridactylus'
$ws.Range("H57").Value = '(Linnaeus, 1758)'
$ws.Range("I57").Value = "'1"
$ws.Range("I57").Style = "Normal"
$ws.Range("K57").Value = "'"
$ws.Range("K57").Style = "Normal"
$ws.Range("L57").Value = 'hane'
$ws.Range("M57").Value = 'födosökande'
$ws.Range("N57").Value = "'"
$ws.Range("N57").Style = "Normal"
$ws.Range("Q57").Value = 468224
$ws.Range("R57").Value = 7017146
$ws.Range("Z57").Value = '14:17'
$ws.Range("AB57").Value = '14:17'
$ws.Range("AC57").Value = 'Bearbetade medelgrov gran. Flög över till klenare gran när jag närmade mig men förvånansvärt obrydd.'

# --- Row 58 ---
$ws.Range("A58").Value = 131030174
$ws.Range("I58").Value = "'"
$ws.Range("I58").Style = "Normal"
$ws.Range("P58").Value = 'Stacktjärnen, Rödön, Stacktjärnen, Rödön, Jmt'
$ws.Range("Q58").Value = 468273
$ws.Range("R58").Value = 7017348
$ws.Range("Z58").Value = '13:30'
$ws.Range("AB58").Value = '13:30'

# --- Row 59 ---
$ws.Range("A59").Value = 131029927
$ws.Range("B59").Value = 58043
$ws.Range("E59").Value = 103021
$ws.Range("F59").Value = 'Talltita'
$ws.Range("G59").Value = 'Poecile montanus'
$ws.Range("H59").Value = '(Conrad von Baldenstein, 1827)'
$ws.Range("I59").Value = "'2"
$ws.Range("I59").Style = "Normal"
$ws.Range("K59").ClearContents()
$ws.Range("L59").ClearContents()
$ws.Range("M59").Value = 'lockläte, övriga läten'
$ws.Range("N59").ClearContents()
$ws.Range("P59").Value = 'Stacktjärnen, Stacktjärnen, Jmt'
$ws.Range("Q59").Value = 468509
$ws.Range("R59").Value = 7017474
$ws.Range("Z59").Value = '13:16'
$ws.Range("AB59").Value = '13:16'
$ws.Range("AC59").ClearContents()
